$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new row 10 data for the "20 bunnies build with 1 triangle per leaf" entry
$ws.Range("B10").Value = "24.8650 seconds"
$ws.Range("C10").Value = 2778041
$ws.Range("D10").Value = 1389021
$ws.Range("F10").Value = 262144
$ws.Range("G10").Value = 11156027
$ws.Range("J10").Value = 42.5569

# Column D width changes from 7 to 8 (stored/"bestFit" width units).
# The ColumnWidth COM property is expressed in character-width units, which
# are offset from the stored XML width by the default column padding
# (~0.8333 chars for this workbook's default font), so 7.17 here round-trips
# to a stored width of 8.
$ws.Columns.Item(4).ColumnWidth = 7.17

# Update the selection to B10
$ws.Range("B10").Select() | Out-Null
